$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.226.50'
$ws.Range('E2').Value = '  +0.06%  '
$ws.Range('D3').Value = '1.904.97'
$ws.Range('E3').Value = '  +0.47%  '
$ws.Range('D4').Value = "'1.002"
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = "'306.30"
$ws.Range('E5').Value = '  -0.34%  '
$ws.Range('D6').Value = "'1.001"
$ws.Range('E6').Value = '  -0.02%  '
$ws.Range('D7').Value = "'0.5364"
$ws.Range('E7').Value = '  +3.16%  '
$ws.Range('D8').Value = "'0.3808"
$ws.Range('D9').Value = "'0.07291"
$ws.Range('E9').Value = '  +0.11%  '
$ws.Range('D10').Value = "'22.25"
$ws.Range('E10').Value = '  +5.05%  '
$ws.Range('D11').Value = "'0.9054"
$ws.Range('E11').Value = '  +0.54%  '
$ws.Range('D12').Value = "'0.08213"
$ws.Range('E12').Value = '  +0.33%  '
$ws.Range('D13').Value = "'95.80"
$ws.Range('E13').Value = '  -0.97%  '
$ws.Range('D14').Value = "'5.343"
$ws.Range('D15').Value = "'1.003"
$ws.Range('E15').Value = '  +0.01%  '
$ws.Range('E16').Value = '  +2.09%  '
$ws.Range('D17').Value = "'0.000008663"
$ws.Range('E17').Value = '  +0.55%  '
$ws.Range('D18').Value = "'1.001"
$ws.Range('E18').Value = '  -0.07%  '
$ws.Range('D19').Value = '27.246.48'
$ws.Range('E19').Value = '  +0.02%  '
$ws.Range('D20').Value = "'5.048"
$ws.Range('E20').Value = '  -0.72%  '
$ws.Range('D21').Value = '1.068.95'
$ws.Range('E21').Value = '  -43.60%  '
$ws.Range('D22').Value = "'10.78"
$ws.Range('E22').Value = '  +0.76%  '
$ws.Range('D23').Value = "'6.520"
$ws.Range('E23').Value = '  +1.83%  '
$ws.Range('D24').Value = "'149.03"
$ws.Range('E24').Value = '  +1.09%  '
$ws.Range('D25').Value = "'2.290"
$ws.Range('E25').Value = '  -0.22%  '
$ws.Range('D26').Value = "'18.39"
$ws.Range('E26').Value = '  +0.89%  '
$ws.Range('D27').Value = "'1.747"
$ws.Range('E27').Value = '  +0.12%  '
$ws.Range('D28').Value = "'116.84"
$ws.Range('E28').Value = '  +1.44%  '
$ws.Range('D29').Value = "'4.816"
$ws.Range('E29').Value = '  -0.26%  '
$ws.Range('D30').Value = "'4.726"
$ws.Range('E30').Value = '  -4.57%  '
$ws.Range('D31').Value = "'0.09224"
$ws.Range('E31').Value = '  -0.03%  '
$ws.Range('D32').Value = "'0.8291"
$ws.Range('E32').Value = '  +4.30%  '
$ws.Range('D33').Value = "'0.05083"
$ws.Range('E33').Value = '  +1.03%  '
$ws.Range('D34').Value = "'1.218"
$ws.Range('E34').Value = '  -0.18%  '
$ws.Range('D35').Value = "'3.006"
$ws.Range('E35').Value = '  +2.10%  '
$ws.Range('D36').Value = "'3.334"
$ws.Range('E36').Value = '  -3.20%  '
$ws.Range('D37').Value = "'2.673"
$ws.Range('E37').Value = '  +3.04%  '
$ws.Range('D38').Value = "'0.5842"
$ws.Range('E38').Value = '  +3.13%  '
$ws.Range('D39').Value = "'0.02002"
$ws.Range('E39').Value = '  +0.81%  '
$ws.Range('E40').Value = '  +0.26%  '
$ws.Range('D41').Value = "'9.312"
$ws.Range('E41').Value = '  +4.03%  '
$ws.Range('D42').Value = "'6.628"
$ws.Range('E42').Value = '  +1.09%  '
$ws.Range('D43').Value = "'117.18"
$ws.Range('E43').Value = '  +1.44%  '
$ws.Range('D44').Value = "'0.5087"
$ws.Range('E44').Value = '  +3.97%  '
$ws.Range('E45').Value = '  +0.45%  '
$ws.Range('D46').Value = "'1.000"
$ws.Range('E46').Value = '  -0.05%  '
$ws.Range('D47').Value = "'10.11"
$ws.Range('E47').Value = '  +0.10%  '
$ws.Range('D48').Value = "'1.641"
$ws.Range('E48').Value = '  +1.05%  '
$ws.Range('D49').Value = "'38.40"
$ws.Range('E49').Value = '  +0.62%  '
$ws.Range('D50').Value = "'0.06157"
$ws.Range('D51').Value = "'63.50"
$ws.Range('E51').Value = '  +0.07%  '
